# Adds a new "2022-Q1" sheet (modeled on the existing "2021-Q4" sheet) right
# before the "总计" (Total) sheet, and updates the "总计" sheet so it gains a
# new leading row for "2022-Q1" (with the rest of its rows shifting down by one).

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (matches the source data, where numeric-looking figures such as "0.46"
# are kept as strings rather than numbers).
function Set-TextValue {
    param($Cell, $Text)
    $helperWs = $wb.Worksheets.Add()
    $helperCell = $helperWs.Range("A1")
    $helperCell.NumberFormat = "@"
    $helperCell.Value = $Text
    $helperCell.Copy()
    $Cell.PasteSpecial(-4163)  # xlPasteValues - copies the value only (keeps it text)
    $helperWs.Delete()
}

# ---------------------------------------------------------------------------
# 1) Remove the existing "总计" sheet so that its sheetId (5) becomes free
#    again; this lets the newly inserted "2022-Q1" sheet take sheetId 5 and
#    the re-created "总计" sheet take sheetId 6, matching the target layout.
# ---------------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

# ---------------------------------------------------------------------------
# 2) Create "2022-Q1" by duplicating "2021-Q4" (same column layout/styles),
#    placed at the end of the workbook, then update the quarter's figures.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$q1new = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1new.Name = "2022-Q1"

# B2 (基金代码/378006) and C2 (基金名称) stay the same as 2021-Q4; update the rest.
Set-TextValue -Cell $q1new.Range("D2") -Text "0.46"
Set-TextValue -Cell $q1new.Range("E2") -Text "88.99"
Set-TextValue -Cell $q1new.Range("F2") -Text "1.79"
Set-TextValue -Cell $q1new.Range("G2") -Text "0.0082"
$q1new.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 3) Re-create "总计" as a fresh sheet right after "2022-Q1" and populate it
#    with the updated rollup (new 2022-Q1 row on top, rest shifted down,
#    2021-Q3's value corrected from 0.02 to 0.01).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1new)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.01

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.01

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.02

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.02

# Apply the bold/centered/bordered header style (matching the other sheets)
# to the header row and the index column (column A) of the "总计" sheet by
# copying formatting from the equivalent cells on "2022-Q1".
$q1new.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats

$q1new.Range("A2").Copy()
$total.Range("A2:A6").PasteSpecial(-4122)  # xlPasteFormats

Write-Host "Workbook updated: added 2022-Q1 sheet and refreshed 总计 sheet."
